# Apply metadata updates to the "Metadata" worksheet of the CodeSystem workbook.
# Changes:
#   - Status:         draft -> active
#   - Experimental:   (blank) -> false   (must stay a literal text value, not a Boolean)
#   - Date:           2025-07-12T17:02:17-03:00 -> 2025-07-14T12:58:17-03:00
#   - Case Sensitive: (blank) -> true    (must stay a literal text value, not a Boolean)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$ws.Range("B6").Value = "active"

# Writing the bare word "false"/"true" via .Value auto-coerces to a Boolean cell,
# so build the literal string through a formula and collapse it back down to a
# plain text value via Copy + PasteSpecial(xlPasteValues) - this keeps the cell's
# original style and produces a genuine text cell instead of a Boolean one.
$expCell = $ws.Range("B7")
$expCell.Formula = "=""false"""
$expCell.Copy()
$expCell.PasteSpecial(-4163)

$ws.Range("B8").Value = "2025-07-14T12:58:17-03:00"

$caseCell = $ws.Range("B15")
$caseCell.Formula = "=""true"""
$caseCell.Copy()
$caseCell.PasteSpecial(-4163)

$excel.CutCopyMode = 0

$wb.Save()
